# =====================================================================
# Update: switzerland_super-league_2023-2024 — re-synced rows from source
# (row order/content changes for 14/15, 43/44/45, 60/61, 64/65, 66/67,
#  72/73, 82/83) and 6 newly scraped fixtures appended as rows 84-89.
# =====================================================================

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column letters F..V, in order, mapped to 1-based column indices 6..22 ---
$cols = @(6,7,8,9,10,11,12,13,14,15,16,17,18,19,20,21,22)

# --- Rows whose F:V contents are replaced in place (A/E — index/date — unchanged) ---
$updates = @(
    @(14, 'Lausanne Ouchy', 1, 'Servette', 1, 3.69, '31/07/2023 05:11', 3.42, '05/08/2023 17:48', 3.91, '31/07/2023 05:11', 3.75, '05/08/2023 17:58', 1.97, '31/07/2023 05:11', 2.1, '05/08/2023 17:53', 'https://www.betexplorer.com/football/switzerland/super-league/lausanne-ouchy-servette/AePC8XxS/'),
    @(15, 'Young Boys', 5, 'Winterthur', 2, 1.26, '30/07/2023 16:42', 1.39, '05/08/2023 17:56', 6.81, '30/07/2023 16:42', 5.57, '05/08/2023 17:57', 10.25, '30/07/2023 16:42', 7.23, '05/08/2023 17:57', 'https://www.betexplorer.com/football/switzerland/super-league/young-boys-winterthur/jk7u4ioq/'),
    @(43, 'Servette', 2, 'Winterthur', 2, 1.56, '24/09/2023 15:42', 1.56, '27/09/2023 20:29', 4.47, '24/09/2023 15:42', 4.44, '27/09/2023 20:29', 4.91, '24/09/2023 15:42', 5.75, '27/09/2023 20:29', 'https://www.betexplorer.com/football/switzerland/super-league/servette-winterthur/vXKqSIZH/'),
    @(44, 'St. Gallen', 2, 'Young Boys', 1, 2.33, '24/09/2023 15:42', 2.38, '27/09/2023 20:29', 3.72, '24/09/2023 15:42', 3.68, '27/09/2023 20:25', 2.75, '24/09/2023 15:42', 2.93, '27/09/2023 20:29', 'https://www.betexplorer.com/football/switzerland/super-league/st-gallen-young-boys/MeAlRxkO/'),
    @(45, 'Lugano', 2, 'Lausanne', 1, 1.81, '24/09/2023 22:12', 2.2, '27/09/2023 19:35', 3.89, '24/09/2023 22:12', 3.84, '27/09/2023 19:35', 3.93, '24/09/2023 22:12', 3.13, '27/09/2023 19:35', 'https://www.betexplorer.com/football/switzerland/super-league/lugano-lausanne/0tMuTbKB/'),
    @(60, 'Lausanne Ouchy', 2, 'Lausanne', 2, 2.88, '11/10/2023 13:43', 2.88, '21/10/2023 17:55', 3.47, '11/10/2023 13:43', 3.48, '21/10/2023 17:55', 2.35, '11/10/2023 13:43', 2.51, '21/10/2023 17:55', 'https://www.betexplorer.com/football/switzerland/super-league/lausanne-ouchy-lausanne/beGieH34/'),
    @(61, 'Young Boys', 0, 'Zurich', 0, 1.76, '11/10/2023 13:43', 2.03, '21/10/2023 17:59', 4.21, '11/10/2023 13:43', 3.7, '21/10/2023 17:59', 3.84, '11/10/2023 13:43', 3.69, '21/10/2023 17:59', 'https://www.betexplorer.com/football/switzerland/super-league/young-boys-zurich/GWAdfyJA/'),
    @(64, 'Luzern', 2, 'Yverdon', 1, 1.58, '11/10/2023 13:43', 1.68, '22/10/2023 16:25', 4.37, '11/10/2023 13:43', 4.33, '22/10/2023 16:25', 4.84, '11/10/2023 13:43', 4.77, '22/10/2023 16:25', 'https://www.betexplorer.com/football/switzerland/super-league/luzern-yverdon/61kymgBp/'),
    @(65, 'Grasshoppers', 2, 'Lugano', 1, 2.64, '11/10/2023 13:43', 2.71, '22/10/2023 16:29', 3.5, '11/10/2023 13:43', 3.7, '22/10/2023 16:22', 2.64, '11/10/2023 13:43', 2.54, '22/10/2023 16:29', 'https://www.betexplorer.com/football/switzerland/super-league/grasshoppers-lugano/0KC8iZ2T/'),
    @(66, 'Yverdon', 1, 'Winterthur', 1, 2.65, '22/10/2023 16:42', 2.47, '28/10/2023 17:59', 3.64, '22/10/2023 16:42', 3.73, '28/10/2023 17:57', 2.44, '22/10/2023 16:42', 2.78, '28/10/2023 17:59', 'https://www.betexplorer.com/football/switzerland/super-league/yverdon-winterthur/8QZZ7fmA/'),
    @(67, 'St. Gallen', 3, 'Grasshoppers', 1, 1.39, '22/10/2023 16:42', 1.62, '28/10/2023 17:36', 5.14, '22/10/2023 16:42', 4.56, '28/10/2023 17:59', 6.33, '22/10/2023 16:42', 5.02, '28/10/2023 17:59', 'https://www.betexplorer.com/football/switzerland/super-league/st-gallen-grasshoppers/fNVV8zY3/'),
    @(72, 'Lausanne', 3, 'Lugano', 1, 2.24, '29/10/2023 16:42', 2.34, '04/11/2023 17:52', 3.67, '29/10/2023 16:42', 3.67, '04/11/2023 17:52', 3.08, '29/10/2023 16:42', 3, '04/11/2023 17:52', 'https://www.betexplorer.com/football/switzerland/super-league/lausanne-lugano/Mkku5hXS/'),
    @(73, 'Winterthur', 1, 'Young Boys', 4, 3.25, '29/10/2023 16:42', 3.47, '04/11/2023 17:57', 4.23, '29/10/2023 16:42', 3.89, '04/11/2023 17:57', 1.93, '29/10/2023 16:42', 2.04, '04/11/2023 17:57', 'https://www.betexplorer.com/football/switzerland/super-league/winterthur-young-boys/WOvFbjAd/'),
    @(82, 'Servette', 4, 'Basel', 1, 1.67, '05/11/2023 16:42', 1.62, '12/11/2023 16:27', 4.26, '05/11/2023 16:42', 4.37, '12/11/2023 16:27', 4.72, '05/11/2023 16:42', 5.27, '12/11/2023 16:25', 'https://www.betexplorer.com/football/switzerland/super-league/servette-basel/KpL9KXWF/'),
    @(83, 'Lugano', 0, 'Zurich', 3, 2.49, '05/11/2023 16:42', 2.93, '12/11/2023 16:23', 3.51, '05/11/2023 16:42', 3.26, '12/11/2023 16:29', 2.81, '05/11/2023 16:42', 2.59, '12/11/2023 16:23', 'https://www.betexplorer.com/football/switzerland/super-league/lugano-zurich/tfM5LDH9/'),
)

foreach ($entry in $updates) {
    $r = $entry[0]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Cells.Item($r, $cols[$i]).Value = $entry[$i + 1]
    }
}

# --- New rows 84-89: copy formatting (borders/font/number-format) from the
#     last existing data row, then populate all columns A:V ---
$lastRow = 83
$firstNewRow = 84
$lastNewRow = 89
$ws.Range("A$lastRow`:V$lastRow").Copy() | Out-Null
$ws.Range("A$firstNewRow`:V$lastNewRow").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$newRows = @(
    @(84, 83, 45255.75, 'Luzern', 3, 'Winterthur', 1, 1.74, '15/11/2023 15:12', 1.87, '25/11/2023 17:59', 4.34, '15/11/2023 15:12', 4.34, '25/11/2023 17:59', 3.69, '15/11/2023 15:12', 3.56, '25/11/2023 17:59', 'https://www.betexplorer.com/football/switzerland/super-league/luzern-winterthur/S8qbbsGD/'),
    @(85, 84, 45255.75, 'Zurich', 3, 'Young Boys', 1, 2.5, '15/11/2023 15:12', 2.5, '25/11/2023 17:59', 3.41, '15/11/2023 15:12', 3.42, '25/11/2023 17:59', 2.66, '15/11/2023 15:12', 2.93, '25/11/2023 17:59', 'https://www.betexplorer.com/football/switzerland/super-league/zurich-young-boys/xvRlL4NQ/'),
    @(86, 85, 45255.85416666666, 'Yverdon', 0, 'Lugano', 5, 2.97, '15/11/2023 15:12', 3.75, '25/11/2023 20:28', 3.5, '15/11/2023 15:12', 3.75, '25/11/2023 20:28', 2.23, '15/11/2023 15:12', 1.99, '25/11/2023 20:24', 'https://www.betexplorer.com/football/switzerland/super-league/yverdon-lugano/fFypMp8K/'),
    @(87, 86, 45256.59375, 'Lausanne', 1, 'Lausanne Ouchy', 0, 1.79, '15/11/2023 15:12', 1.8, '26/11/2023 13:22', 3.84, '15/11/2023 15:12', 3.99, '26/11/2023 13:54', 3.89, '15/11/2023 15:12', 4.36, '26/11/2023 13:22', 'https://www.betexplorer.com/football/switzerland/super-league/lausanne-lausanne-ouchy/lIpfaN07/'),
    @(88, 87, 45256.6875, 'Servette', 2, 'Grasshoppers', 0, 1.52, '15/11/2023 15:12', 1.56, '26/11/2023 16:25', 4.42, '15/11/2023 15:12', 4.35, '26/11/2023 16:25', 5.07, '15/11/2023 15:12', 5.57, '26/11/2023 16:25', 'https://www.betexplorer.com/football/switzerland/super-league/servette-grasshoppers/tCXuNQhE/'),
    @(89, 88, 45256.6875, 'Basel', 2, 'St. Gallen', 0, 3.51, '15/11/2023 15:12', 3.19, '26/11/2023 16:29', 3.99, '15/11/2023 15:12', 3.85, '26/11/2023 16:29', 1.85, '15/11/2023 15:12', 2.17, '26/11/2023 16:29', 'https://www.betexplorer.com/football/switzerland/super-league/basel-st-gallen/Wtjk03o1/'),
)

foreach ($entry in $newRows) {
    $r = $entry[0]
    $ws.Cells.Item($r, 1).Value = $entry[1]
    $ws.Cells.Item($r, 2).Value = 'switzerland'
    $ws.Cells.Item($r, 3).Value = 'super-league'
    $ws.Cells.Item($r, 4).Value = '2023-2024'
    $ws.Cells.Item($r, 5).Value = $entry[2]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Cells.Item($r, $cols[$i]).Value = $entry[$i + 3]
    }
}

# --- Sheet dimension now spans through the new last row/col ---
$ws.Range("A1:V89").Select()

